## edit.ps1 - applies the customer-letter edits described by the diff:
##   1. Split the "LATORRE 98 / 2290000 CALERA" address line into two
##      lines ("LATORRE 98" / "2290000 CALERA") joined by a line break.
##   2. Resize the three columns of the first (Product Id / Description /
##      Serial-Lot) table from 3116/3117/3117 dxa to 2245/4680/2425 dxa.

$d = $word.ActiveDocument

# --- 1. Split the recipient address line with a manual line break -------
# Word represents a manual line break inside a run's text with Chr(11)
# (vertical-tab). Replacing the matched range's .Text (rather than doing
# a Find/Replace across the whole story) keeps this edit scoped to just
# the "LATORRE..." <w:t>, leaving the neighbouring "Carlos Barroso" /
# "SERVICLINICA S.A." runs untouched.
$addrRange = $d.Content.Duplicate
$addrRange.Find.Execute("LATORRE 98 / 2290000 CALERA", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($addrRange.Find.Found) {
    $target = $d.Range($addrRange.Start, $addrRange.End)
    $lineBreak = [char]11
    $target.Text = "LATORRE 98" + $lineBreak + "2290000 CALERA"
}

# --- 2. Resize the Product Id / Product Description / Serial-Lot table --
$productTable = $d.Tables.Item(1)
$productTable.Columns.Item(1).Width = 112.25   # 2245 dxa
$productTable.Columns.Item(2).Width = 234.0    # 4680 dxa
$productTable.Columns.Item(3).Width = 121.25   # 2425 dxa
